$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 66.36304333333334
$ws.Range("H2").Value = 199.08913
$ws.Range("I2").Value = 0.1775372952319303
$ws.Range("J2").Value = 0.1775372952319303
$ws.Range("M2").Value = 0.5273163333333333
$ws.Range("N2").Value = 1.581949
$ws.Range("Q2").Value = 34.99431667937444
$ws.Range("R2").Value = 314.94885011437
$ws.Range("S2").Value = 0.1775372952319303
$ws.Range("T2").Value = 0.1775372952319303

$ws.Range("G3").Value = 296.1091513333333
$ws.Range("H3").Value = 888.327454
$ws.Range("I3").Value = 0.7921640597024409
$ws.Range("J3").Value = 0.7921640597024407
$ws.Range("M3").Value = 0.5273163333333333
$ws.Range("N3").Value = 1.581949
$ws.Range("Q3").Value = 156.1431919475384
$ws.Range("R3").Value = 1405.288727527846
$ws.Range("S3").Value = 0.7921640597024409
$ws.Range("T3").Value = 0.7921640597024407

$ws.Range("G4").Value = 11.32556566666667
$ws.Range("H4").Value = 33.976697
$ws.Range("I4").Value = 0.03029864506562886
$ws.Range("J4").Value = 0.03029864506562885
$ws.Range("M4").Value = 0.5273163333333333
$ws.Range("N4").Value = 1.581949
$ws.Range("Q4").Value = 5.972155760272556
$ws.Range("R4").Value = 53.749401842453
$ws.Range("S4").Value = 0.03029864506562886
$ws.Range("T4").Value = 0.03029864506562885
